# Requirements use cases gemaakt
#
# Slide 1 (sldId 256) / Slide 2 (sldId 257): the nav-bar textbox's last run
# "<tab><tab> login | register | profile" is split into two separate edits
# (undo of an earlier combined edit):
#   - Slide 1, shape id 9  (TextBox 8) -> "<tab><tab><tab> login | register"
#   - Slide 2, shape id 10 (TextBox 9) -> "<tab><tab> <tab><tab>profile"

function Find-ShapeById($slide, $id) {
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $sh = $slide.Shapes.Item($i)
        if ($sh.Id -eq $id) {
            return $sh
        }
    }
    return $null
}

function Update-NavText($slide, $shapeId, $newTail) {
    $sh = Find-ShapeById $slide $shapeId
    $tr = $sh.TextFrame.TextRange
    $full = $tr.Text
    $oldTail = "`t`t login | register | profile"
    $idx = $full.IndexOf($oldTail)
    if ($idx -ge 0) {
        $sub = $tr.Characters($idx + 1, $oldTail.Length)
        $sub.Text = $newTail
    }
}

$p = $ppt.ActivePresentation

$slide256 = $p.Slides.Item(1)
$slide257 = $p.Slides.Item(2)

Update-NavText $slide256 9 "`t`t`t login | register"
Update-NavText $slide257 10 "`t`t `t`tprofile"
